$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "order" and "IsFood" columns (C and D) are no longer needed - drop them.
# Deleting the columns shifts the remaining data left and automatically
# shrinks the sheet dimension / row spans / prunes the shared string table.
$ws.Columns("C:D").Delete() | Out-Null

# Replace the remaining two-column table (id / Name) with the new cargo list:
# energy/iron/waterdrop stay, the military units (tank/chopper/ship) are
# replaced with fighter/bomber/laser, and waterdrop's Chinese name is
# shortened while a new "antimatter" entry is introduced for energy.
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "Name"

$ws.Range("A2").Value = "energy"
$ws.Range("B2").Value = "反物质"

$ws.Range("A3").Value = "iron"
$ws.Range("B3").Value = "铁"

$ws.Range("A4").Value = "fighter"
$ws.Range("B4").Value = "战斗机"

$ws.Range("A5").Value = "bomber"
$ws.Range("B5").Value = "轰炸机"

$ws.Range("A6").Value = "laser"
$ws.Range("B6").Value = "激光炮"

$ws.Range("A7").Value = "waterdrop"
$ws.Range("B7").Value = "水滴"

# Rebuild the AutoFilter so it only spans the remaining A:B columns (and no
# longer carries the old per-filter sortState that referenced column D).
$ws.AutoFilterMode = $false
$ws.Range("A1:B1").AutoFilter() | Out-Null

# Keep the workbook-level hidden _FilterDatabase name in sync with the
# shrunk AutoFilter range.
$name = $wb.Names.Item(1)
$name.RefersTo = "=工作表1!`$A`$1:`$B`$1"

# Match the saved selection/active cell.
$ws.Range("B3").Select() | Out-Null
